$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Attendance" column (column AB) entirely - the cells to its
# right (Faculty_Feedback, English_Read) shift left to fill the gap.
$ws.Range("AB1").EntireColumn.Delete()
